# Applies the update described by the commit "Atualizado por script em 29-11-2023 14:45"
# - Swaps the match data (columns F:V) between rows 154 and 156
# - Swaps the match data (columns F:V) between rows 155 and 157
# - Appends a new match row (159): Marek vs Chernomorets 1919

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap F:V content between row 154 and row 156 ---
$row154 = $ws.Range("F154:V154").Value2
$row156 = $ws.Range("F156:V156").Value2
$ws.Range("F154:V154").Value2 = $row156
$ws.Range("F156:V156").Value2 = $row154

# --- Swap F:V content between row 155 and row 157 ---
$row155 = $ws.Range("F155:V155").Value2
$row157 = $ws.Range("F157:V157").Value2
$ws.Range("F155:V155").Value2 = $row157
$ws.Range("F157:V157").Value2 = $row155

# --- Append new row 159 with the same look & feel as the previous data row (158) ---
$ws.Cells.Item(158,1).Copy()
$ws.Cells.Item(159,1).PasteSpecial(-4122)
$ws.Cells.Item(158,5).Copy()
$ws.Cells.Item(159,5).PasteSpecial(-4122)

$ws.Cells.Item(159,1).Value2 = 158
$ws.Cells.Item(159,2).Value2 = "bulgaria"
$ws.Cells.Item(159,3).Value2 = "vtora-liga"
$ws.Cells.Item(159,4).Value2 = "2023-2024"
$ws.Cells.Item(159,5).Value2 = 45259.5625
$ws.Cells.Item(159,6).Value2 = "Marek"
$ws.Cells.Item(159,7).Value2 = 3
$ws.Cells.Item(159,8).Value2 = "Chernomorets 1919"
$ws.Cells.Item(159,9).Value2 = 1
$ws.Cells.Item(159,10).Value2 = 1.95
$ws.Cells.Item(159,11).Value2 = "28/11/2023 02:42"
$ws.Cells.Item(159,12).Value2 = 1.85
$ws.Cells.Item(159,13).Value2 = "29/11/2023 13:00"
$ws.Cells.Item(159,14).Value2 = 3.07
$ws.Cells.Item(159,15).Value2 = "28/11/2023 02:42"
$ws.Cells.Item(159,16).Value2 = 3.01
$ws.Cells.Item(159,17).Value2 = "29/11/2023 13:20"
$ws.Cells.Item(159,18).Value2 = 3.54
$ws.Cells.Item(159,19).Value2 = "28/11/2023 02:42"
$ws.Cells.Item(159,20).Value2 = 4.3
$ws.Cells.Item(159,21).Value2 = "29/11/2023 13:00"
$ws.Cells.Item(159,22).Value2 = "https://www.betexplorer.com/football/bulgaria/vtora-liga/marek-chernomorets-1919/Wr9L8o1C/"
